$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table shrinks from 4 data rows to 1 data row. The surviving row keeps
# the book that used to be on row 3 (2663854991 / "Quia aut autem."), so
# deleting the original row 2 shifts row 3's data (and its string cell
# types) up into row 2 rather than re-typing the values by hand.
$ws.Rows.Item(2).Delete()

# Its borrow_count is updated from 3 to 1.
$ws.Range("C2").Value = 1

# Drop the remaining now-stale rows (originally rows 4 and 5, now shifted
# up to rows 3 and 3 again after each delete) so only header + 1 row remain.
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()
